$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")
$ws.Activate()

# 1) Fix the "Note: The placeholders ..." cell (A3): the original had two
#    differently-formatted text runs; collapse into one plain string, fix
#    wording ("or BT background" -> ", BT background") and add trailing \n
#    like the other "tn:Note:" rows in this sheet.
$ws.Cells.Item(3, 1).Value = "tn:Note: The placeholders {ET}, {BT}, {time}, {ETB}, {BTB} will be substituted by the current ET, BT, time, ET background, BT background value in Serial/CallProgram/MODBUS/S7/WebSocket commands`n"
$ws.Rows.Item(3).RowHeight = 13.8

# 2) Insert a new row for the PHIDGET DCMotor "limit" IO command right after
#    the existing "vel(c,v[,sn])" row (row 29), pushing everything below it
#    down by one row.
$ws.Rows.Item(30).Insert()
$ws.Cells.Item(30, 2).Value = "limit(c,v[,sn])"
$ws.Cells.Item(30, 3).Value = "PHIDGET DCMotor: sets current limit of channel c to v (float)"
$ws.Rows.Item(30).RowHeight = 13.8
